$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.597.93'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '1.585.97'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.81'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.42%  '
$ws.Range("E6").Value = '  -2.64%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -2.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0616'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.58'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.08%  '
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").Value = '1.808.32'
$ws.Range("E12").Value = '  -2.72%  '
$ws.Range("D13").Value = '1.586.16'
$ws.Range("E13").Value = '  -2.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.526'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.75'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '26.588.67'
$ws.Range("E17").Value = '  -2.21%  '
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '206.94'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.73'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("E22").Value = '  -3.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.37'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.87'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.28'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  +0.78%  '
$ws.Range("E28").Value = '  -2.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.26'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0506'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("E31").Value = '  -2.09%  '
$ws.Range("E32").Value = '  -4.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.665'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +23.33%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.92'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.87%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '1.327.67'
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.50'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.20%  '
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("E38").Value = '  -1.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.825'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.35'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.783'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.09%  '
$ws.Range("E43").Value = '  -3.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.43'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").Value = '1.721.46'
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.87'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.830'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.21%  '
$ws.Range("E49").Value = '  -1.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0980'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.51'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.71%  '
